$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column G (direccion_origen) before the current comuna column.
$ws.Range("G1").EntireColumn.Insert()

# Insert two new columns after url_1 (currently K, will be L after the insert above)
# for similitud_nombre_dominio and comparacion_direccion.
$ws.Range("L1:M1").EntireColumn.Insert()

# Header row
$ws.Range("G1").Value = "direccion_origen"
$ws.Range("L1").Value = "similitud_nombre_dominio"
$ws.Range("M1").Value = "comparacion_direccion"

# Update row 2 (ABSIDE S.A.) data values
$ws.Range("C2").Value = "+56979776350, +56912345678"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "https://abside.com"
$ws.Range("F2").Value = "Avenida Einstein 863"
$ws.Range("G2").Value = "av. einstein 863"
$ws.Range("H2").Value = "RECOLETA"
$ws.Range("I2").Value = "XIII REGION METROPOLITANA"
$ws.Range("J2").Value = "Desde 1988 tu partner en la construcción. Somos una empresa innovadora que ofrece soluciones de serie y a la medida para todo tipo de proyectos de construcción. Actividades Económicas o Giros VENTA AL POR MAYOR DE MAQUINARIA PARA LA CONSTRUCCION. ARRIENDO DE MAQUINARIA PARA CONSTRUCCION."
$ws.Range("K2").Value = "https://www.genealog.cl/Geneanexus/empresa/CHILE/TNzk3TwNzYzNTAtMA-jTw/nombre-y-rut/ABSIDE-S.A.-79776350-0"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "76.9%"
$ws.Range("M2").Value = "Exacto"
$ws.Range("N2").Value = ""

# Remove row 3 (ACEVEDO Y CIA. LTDA.) entirely
$ws.Range("A3:N3").EntireRow.Delete()
